# Append the Nifty Smallcap 50 data update (through 2025-01-14): 73 new
# trading-day rows (3409-3481) with Date / Total Returns Index / P/E / P/B /
# Div Yield, continuing directly after the existing last row (3408, date
# 45562 = 2024-09-27).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(45565, 11408.71, 32.43, 3.74, 0.8),
  @(45566, 11509.28, 32.72, 3.77, 0.8),
  @(45568, 11288.76, 32.09, 3.7, 0.8100000000000001),
  @(45569, 11198.88, 31.84, 3.67, 0.82),
  @(45572, 10910.09, 31.02, 3.57, 0.84),
  @(45573, 11093.36, 31.54, 3.63, 0.83),
  @(45574, 11235.25, 31.94, 3.68, 0.82),
  @(45575, 11269.31, 32.97, 3.87, 0.8),
  @(45576, 11356.06, 33.23, 3.9, 0.79),
  @(45579, 11430.37, 33.44, 3.92, 0.79),
  @(45580, 11539.41, 33.76, 3.96, 0.78),
  @(45581, 11532.27, 33.74, 3.96, 0.78),
  @(45582, 11363.98, 33.33, 3.9, 0.79),
  @(45583, 11372.59, 33.35, 3.9, 0.79),
  @(45586, 11177.01, 32.64, 3.84, 0.79),
  @(45587, 10763.82, 30.99, 3.69, 0.82),
  @(45588, 10889.12, 31.31, 3.74, 0.8100000000000001),
  @(45589, 10898.6, 31.32, 3.74, 0.8100000000000001),
  @(45590, 10666.44, 30.59, 3.66, 0.83),
  @(45593, 10811.11, 31.41, 3.71, 0.82),
  @(45594, 10912.54, 31.6, 3.75, 0.8100000000000001),
  @(45595, 10999.17, 31.85, 3.78, 0.8),
  @(45596, 11135.82, 32.16, 3.82, 0.8),
  @(45597, 11261.59, 32.52, 3.87, 0.78),
  @(45600, 11032.36, 31.88, 3.79, 0.79),
  @(45601, 11097.86, 32.08, 3.81, 0.79),
  @(45602, 11322.85, 32.72, 3.89, 0.79),
  @(45603, 11196.34, 32.34, 3.84, 0.79),
  @(45604, 11057.25, 31.93, 3.79, 0.82),
  @(45607, 10924.11, 31.56, 3.75, 0.8100000000000001),
  @(45608, 10755.54, 31.13, 3.69, 0.82),
  @(45609, 10449.08, 30.26, 3.58, 0.84),
  @(45610, 10537.83, 30.32, 3.61, 0.85),
  @(45614, 10476.3, 30.11, 3.59, 0.85),
  @(45615, 10586.53, 29.48, 3.63, 0.85),
  @(45617, 10558.49, 29.37, 3.62, 0.86),
  @(45618, 10652.99, 29.63, 3.65, 0.84),
  @(45621, 10817.43, 30.09, 3.71, 0.8100000000000001),
  @(45622, 10923.16, 30.38, 3.75, 0.8100000000000001),
  @(45623, 11052.22, 30.74, 3.75, 0.8),
  @(45624, 11038.19, 30.7, 3.72, 0.8),
  @(45625, 11148.98, 31.12, 3.78, 0.84),
  @(45628, 11280.06, 31.49, 3.82, 0.83),
  @(45629, 11384.07, 31.78, 3.86, 0.82),
  @(45630, 11450.29, 31.96, 3.88, 0.82),
  @(45631, 11601.62, 32.39, 3.93, 0.8100000000000001),
  @(45632, 11698.55, 32.66, 3.96, 0.8),
  @(45635, 11710.49, 32.69, 3.97, 0.8),
  @(45636, 11729.55, 32.74, 3.98, 0.8),
  @(45637, 11747.39, 32.79, 3.98, 0.8),
  @(45638, 11613.8, 32.42, 3.94, 0.8100000000000001),
  @(45639, 11592.66, 32.36, 3.93, 0.8100000000000001),
  @(45642, 11637.3, 32.49, 3.94, 0.8100000000000001),
  @(45643, 11539.42, 32.21, 3.91, 0.8100000000000001),
  @(45644, 11437.13, 31.93, 3.88, 0.82),
  @(45645, 11389.11, 31.79, 3.86, 0.82),
  @(45646, 11075.66, 30.92, 3.75, 0.85),
  @(45649, 11035.82, 30.81, 3.74, 0.85),
  @(45650, 11085.87, 30.95, 3.76, 0.85),
  @(45652, 11079.78, 30.93, 3.75, 0.85),
  @(45653, 11096.71, 30.98, 3.76, 0.85),
  @(45656, 11007.22, 30.73, 3.73, 0.85),
  @(45657, 11089.76, 30.88, 3.75, 0.84),
  @(45658, 11177.98, 31.12, 3.78, 0.84),
  @(45659, 11248.56, 31.32, 3.8, 0.83),
  @(45660, 11183.46, 31.14, 3.78, 0.83),
  @(45663, 10840.58, 30.18, 3.66, 0.86),
  @(45664, 10967.71, 30.54, 3.71, 0.85),
  @(45665, 10798.55, 30.07, 3.65, 0.86),
  @(45666, 10681.45, 29.74, 3.61, 0.87),
  @(45667, 10406.99, 28.98, 3.52, 0.9),
  @(45670, 10011.45, 27.88, 3.38, 0.93),
  @(45671, 10211.32, 28.44, 3.45, 0.91)
)

$startRow = 3409
for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = $startRow + $i
  $row = $rows[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
}

Write-Output ("Rows count: " + $rows.Length)
Write-Output ("Last row written: " + ($startRow + $rows.Length - 1))